# "Generate Report for Archive"
# The localization status for the tracked file moved from "Ready for
# handoff" to "In Translation". That text lives in three places:
#   - Overview!E2 and Overview!F2 (per-language status columns)
#   - zh-cn!C2  and  de-de!C2     (the "Status" column of each lang sheet)
# Excel auto-shrinks the "Status"/lang columns afterwards because the new
# text is shorter than the old, so we re-fit those columns too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Re-fit the columns that held the now-shorter status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
